# "Added ability to set sides of screen as player." (Closes #34)
# The Time Log workbook gets a new logged time-tracking entry on Sheet1
# (row 103): 2014-10-30, 19:30 -> 21:07, 0 min interruption, Category =
# Coding. Downstream formulas (Sheet1!E122 total, Sheet2 SUMIF/percentage
# table, and the pie chart that reads from Sheet2) recalculate from this.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of data, matching the columns/styles used by the surrounding rows.
$ws.Range("A103").Value = 41942                     # Date: 2014-10-30
$ws.Range("B103").Value = 0.8125                     # Start Time: 7:30 PM
$ws.Range("C103").Value = 0.87986111111111109        # Stop Time: 9:07 PM
$ws.Range("D103").Value = 0                          # Interruption: 0 mins
$ws.Range("E103").Formula = "=IF(AND(NOT(ISBLANK(B103)),NOT(ISBLANK(C103))), (C103-B103) * 24 - D103/60, """")"
$ws.Range("F103").Value = "Coding"                   # Activity/Category

# Recalculate so the Sheet1 total, Sheet2 summary table and chart all
# pick up the new row.
$excel.CalculateFullRebuild()

# Move the visible window / selection the way the author left it (scrolled
# up a bit, with A104 selected instead of B121).
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 93
$win.ScrollColumn = 1
$ws.Range("A104").Select() | Out-Null
